$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SP")

# Two new Stored Procedure rows appended after the existing last row (50),
# following the same A/B/C (Name / Parameter / Comment) layout used by the
# rest of the "SP" sheet.
$ws.Range("A51").Value = "Usp_Cp_ForeignKeyControl_Upd"
$ws.Range("B51").Value = "int TBSDYF, String empNo,int Switch"
$ws.Range("C51").Value = "(每日複製)控制外來鍵"

$ws.Range("A52").Value = "Usp_Cp_CdCode_Ins"
$ws.Range("B52").Value = "String EmpNo"
$ws.Range("C52").Value = "(每日複製)CdCode"

# The sheet's selection moves down by the same two rows that were inserted.
$ws.Range("C53").Select() | Out-Null
